$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "제보) 대학교의 어설픈 인공지능/데이터 사이언스 교육이 학생들을 망치고 있습니다"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/university-failure-in-data-science-education/#utm_source=rss&utm_medium=rss&utm_campaign=university-failure-in-data-science-education"

$ws.Range("D28").Value = "ADAM : A METHOD FOR STOCHASTIC OPTIMIZATION 리뷰 (작성중)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/90"

$ws.Range("D37").Value = "[Paper Review] OOD Detection Using an Ensemble of Self-Supervised Leave-out Classifiers"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1442&mod=document&pageid=1"

$ws.Range("D39").Value = "Autoencoder: Neural Networks For Unsupervised Learning"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Autoencoder-Neural-Networks-For-Unsupervised-Learning"

$ws.Range("D46").Value = "Chest X-Ray Medical Diagnosis with Deep Learning - ⑤ Prediction and Evaluation"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/385"

$ws.Range("D51").Value = "[세이버메트릭스] 타타타자로 시작하는 말- 타석, 타수, 타율"
$ws.Range("E51").Value = "https://bskyvision.com/1010"
